$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A/B columns for the new rows (ID + Page), rows 26-31 ---
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Interactive Binning"
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Interactive Binning"
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Interactive Binning"
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "Interactive Binning"
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "Interactive Binning"
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Interactive Binning"

# --- Feature descriptions typed first, top to bottom (column C) ---
$ws.Range("C25").Value = "Show a dropdown for selecting an automated binning algorithm"
$ws.Range("C26").Value = "Update the description of the automated binning algorithm selected based on dropdown value"
$ws.Range("C27").Value = "Show mixed chart"
$ws.Range("C28").Value = "Show control panel"
$ws.Range("C29").Value = "Show initial (after) summary statistics table panel"

# --- Then the rest of row 25's fields ---
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "Closed"
$ws.Range("F25").Value = 44978
$ws.Range("H25").Value = "-"

# --- Then the rest of row 26's fields ---
$ws.Range("D26").Value = "-"
$ws.Range("E26").Value = "Closed"
$ws.Range("F26").Value = 44978
$ws.Range("H26").Value = "-"

# --- Then the rest of row 27's fields ---
$ws.Range("D27").Value = "-"
$ws.Range("E27").Value = "In Progress"

# --- Then the rest of row 28's fields ---
$ws.Range("D28").Value = "-"
$ws.Range("E28").Value = "In Progress"

# --- Then the rest of row 29's fields ---
$ws.Range("D29").Value = "-"

# --- Selection moved to C32 (where the author was about to type next) ---
$ws.Range("C32").Select()
